$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellXml($row, $col, $bodyXml) {
    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

$body_2_1 = '<w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ϣⲟⲙⲧ ⲛ̀ⲣⲁⲛ ⲉⲧϧⲉⲛ ⲛⲓⲫⲏⲟⲩⲓ̀:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲛ̀ⲑⲟⲕ ⲁⲕⲉⲣⲫⲟⲣⲓⲛ ⲙ̀ⲙⲱⲟⲩ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲡⲓⲑⲉⲱⲣⲓⲙⲟⲥ ⲛ̀ⲉⲩⲁⲅⲅⲉⲗⲓⲥⲧⲏⲥ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ⲙⲁⲣⲕⲟⲥ ⲡⲓⲁ̀ⲡⲟⲥⲧⲟⲗⲟⲥ.</w:t></w:r></w:p>'
Set-CellXml 2 1 $body_2_1

$body_2_2 = '<w:p><w:r><w:t>You wore three names:</w:t></w:r></w:p><w:p><w:r><w:t>That are in heaven:</w:t></w:r></w:p><w:p><w:r><w:t>O Behold of God:</w:t></w:r></w:p><w:p><w:r><w:t>The Evangelist:</w:t></w:r></w:p><w:p><w:r><w:t>Mark the Apostle.</w:t></w:r></w:p>'
Set-CellXml 2 2 $body_2_2

$body_2_3 = '<w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>You bear three names</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>In heaven,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>O Beholder of God, the Evangelist,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>Mark the Apostle.</w:t></w:r></w:p>'
Set-CellXml 2 3 $body_2_3

$body_3_1 = '<w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ⲁⲕⲉⲣⲫⲟⲣⲓⲛ ⲙ̀ⲡⲓϣⲟⲙⲧ ⲛ̀ⲭ̀ⲗⲟⲙ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲡⲓϣⲟⲙⲧ ⲛ̀ⲣⲁⲛ ⲉⲧϫⲏⲕ ⲉ̀ⲃⲟⲗ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲉ̀ⲧⲉ ⲫⲁⲓ ⲡⲉ Ⲫⲓⲱⲧ ⲛⲉ Ⲡϣⲏⲉⲣⲓ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲛⲉⲙ Ⲡⲓⲡ̀ⲛⲉⲩⲙⲁ ⲉⲑⲟⲩⲁⲃ.</w:t></w:r></w:p>'
Set-CellXml 3 1 $body_3_1

$body_3_2 = '<w:p><w:r><w:t>You wore the three crowns:</w:t></w:r></w:p><w:p><w:r><w:t>The three perfect names:</w:t></w:r></w:p><w:p><w:r><w:t>Which is the Father, the Son:</w:t></w:r></w:p><w:p><w:r><w:t>And the Holy Spirit.</w:t></w:r></w:p>'
Set-CellXml 3 2 $body_3_2

$body_3_3 = '<w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>You wear three crowns,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>In the perfect Name of</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>The Father, the Son</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>And the Holy Spirit:</w:t></w:r></w:p>'
Set-CellXml 3 3 $body_3_3

$body_4_1 = '<w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ⲛⲑⲟⲕ ⲟⲩⲁ̀ⲡⲟⲥⲧⲟⲗⲟⲥ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲛ̀ⲑⲟⲕ ⲟⲛ ⲟⲩⲙⲁⲣⲧⲩⲣⲟⲥ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲛ̀ⲑⲟⲕ ⲟⲛ ⲡⲉ ⲛⲓⲙⲁϩ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲥ̀ⲛⲁⲩ ⲛ̀ⲥⲱⲧⲡ ⲛ̀ⲉⲩⲁⲅⲅⲉⲗⲓⲥⲧⲏⲥ.</w:t></w:r></w:p>'
Set-CellXml 4 1 $body_4_1

$body_4_2 = '<w:p><w:r><w:t>You are an apostle:</w:t></w:r></w:p><w:p><w:r><w:t>You are a martyr:</w:t></w:r></w:p><w:p><w:r><w:t>You are also the second</w:t></w:r></w:p><w:p><w:r><w:t>Chosen Evangelist.</w:t></w:r></w:p>'
Set-CellXml 4 2 $body_4_2

$body_4_3 = '<w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>You are an Apostle,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>You are a martyr,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>And you are the second</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>Elect Evangelist.</w:t></w:r></w:p>'
Set-CellXml 4 3 $body_4_3

$body_5_1 = '<w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ⲛⲉⲕⲕⲉϣ̀ⲫⲏⲣ ⲛ̀ⲁ̀ⲡⲟⲥⲧⲟⲗⲟⲥ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t xml:space="preserve">ⲥⲉϣⲟⲩϣⲟⲩ ⲙ̀ⲙⲱⲟⲩ ⲉ̀ϩ̀ⲣⲏⲓ </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>ⲉ̀ϫⲱⲕ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲟⲩⲟϩ ⲛⲉⲕⲥⲁϫⲓ ⲁⲩⲫⲟϩ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ϣⲁ ⲁⲩⲫⲣⲏϫⲥ ⲛ̀ϯⲟⲓⲕⲟⲙⲉⲛⲏ.</w:t></w:r></w:p>'
Set-CellXml 5 1 $body_5_1

$body_5_2 = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Your Apostolic friends:</w:t></w:r></w:p><w:p><w:r><w:t>Boast about you:</w:t></w:r></w:p><w:p><w:r><w:t>And your words reached:</w:t></w:r></w:p><w:p><w:r><w:t>The ends of the world.</w:t></w:r></w:p>'
Set-CellXml 5 2 $body_5_2

$body_5_3 = '<w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>Your Apostolic compatriots</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Boast of you,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>For your words have reached</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>The ends of the world.</w:t></w:r></w:p>'
Set-CellXml 5 3 $body_5_3

$body_6_1 = '<w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Ⲥⲉϣⲟⲩϣⲟⲩ ⲙ̀ⲙⲱⲟⲩ ⲛ̀ϧ̀ⲣⲏⲓ ⲛ̀ϭⲏⲧⲕ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲛ̀ϫⲉ ⲛⲏⲉⲧⲁⲕⲧⲟϫⲟⲩ ϩⲓϫⲉⲛ ⲡⲓⲕⲁϩⲓ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ϧⲉⲛ ⲧ̀ⲭⲱⲣⲁ ⲧⲏⲣⲥ ⲛ̀Ⲭⲏⲙⲓ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲁⲩⲫⲓⲣⲓ ⲉ̀ⲃⲟⲗ ⲉⲩϯⲕⲁⲣⲡⲟⲥ.</w:t></w:r></w:p>'
Set-CellXml 6 1 $body_6_1

$body_6_2 = '<w:p><w:r><w:t>Through whom you planted on earth:</w:t></w:r></w:p><w:p><w:r><w:t>Through the land of Egypt:</w:t></w:r></w:p><w:p><w:r><w:t>Take pride in you:</w:t></w:r></w:p><w:p><w:r><w:t>They blossomed and brought forth fruit.</w:t></w:r></w:p>'
Set-CellXml 6 2 $body_6_2

$body_6_3 = '<w:p><w:pPr><w:pStyle w:val="EngHangEnd"/><w:ind w:left="0" w:firstLine="0"/></w:pPr><w:r><w:t>Y</w:t></w:r><w:r><w:t xml:space="preserve">ou planted </w:t></w:r><w:r><w:t xml:space="preserve">these words </w:t></w:r><w:r><w:t>on earth,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>And throughout Egypt,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>(Which takes pride in you),</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>And they blossomed and bore fruit.</w:t></w:r></w:p>'
Set-CellXml 6 3 $body_6_3

$body_7_1 = '<w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ⲧⲱⲃϩ ⲙ̀Ⲡⲟ̄ⲥ̄ ⲉ̀ϩ̀ⲣⲏⲓ ⲉ̀ϫⲱⲛ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲱ̀ ⲡⲓⲑⲉⲱⲣⲓⲙⲟⲥ ⲛ̀ⲉⲩⲁⲅⲅⲉⲗⲓⲥⲧⲏⲥ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>Ⲙⲁⲣⲕⲟⲥ ⲡⲓⲁ̀ⲡⲟⲥⲧⲟⲗⲟⲥ:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="CopticHangingVerse"/></w:pPr><w:r><w:t>ⲛ̀ⲧⲉⲫⲭⲁ ⲛⲉⲛⲛⲟⲃⲓ ⲛⲁⲛ ⲉ̀ⲃⲟⲗ.</w:t></w:r></w:p>'
Set-CellXml 7 1 $body_7_1

$body_7_2 = '<w:p><w:r><w:t>Pray to the Lord on our behalf:</w:t></w:r></w:p><w:p><w:r><w:t>O Beholder of God, the Evangelist:</w:t></w:r></w:p><w:p><w:r><w:t>Mark the apostle:</w:t></w:r></w:p><w:p><w:r><w:t>That He may forgive us our sins.</w:t></w:r></w:p>'
Set-CellXml 7 2 $body_7_2

$body_7_3 = '<w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>Pray to the Lord on our behalf,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>O Beholder of God, the Evangelist,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>Mark the Apostle,</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="EngHangEnd"/></w:pPr><w:r><w:t>That He may forgive us our sins.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Set-CellXml 7 3 $body_7_3
